# errori_lettura.xlsx -- rebuild rows 2-82 ("commessa" list) to match the new export.
# The edit is a full reshuffle of the 81 data rows (same set of "commessa" ids,
# new row order) plus a handful of per-row flag/highlight tweaks.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: fix up the yellow "Interior" highlight on column A before the value rewrite ---
# (done first, using cells whose CURRENT highlight state is still correct, as format donors,
#  so every other row keeps reusing the workbook existing style index unchanged)
$ws.Range("A62").Copy()
$ws.Range("A56").PasteSpecial(-4122)
$ws.Range("A57").PasteSpecial(-4122)
$ws.Range("A61").PasteSpecial(-4122)

$ws.Range("A2").Copy()
$ws.Range("A62").PasteSpecial(-4122)
$ws.Range("A63").PasteSpecial(-4122)

# --- Step 2: rewrite every data row (2-82): "commessa" id in column A plus the
#     "CAMPO VUOTO" marker cells in H / L / O / P for that row ---
$ws.Range("A2").Value = 245089
$ws.Range("H2").Value = "CAMPO VUOTO"
$ws.Range("L2").ClearContents()
$ws.Range("O2").Value = "CAMPO VUOTO"
$ws.Range("P2").Value = "CAMPO VUOTO"

$ws.Range("A3").Value = 252713
$ws.Range("H3").ClearContents()
$ws.Range("L3").ClearContents()
$ws.Range("O3").Value = "CAMPO VUOTO"
$ws.Range("P3").Value = "CAMPO VUOTO"

$ws.Range("A4").Value = 251632
$ws.Range("H4").ClearContents()
$ws.Range("L4").ClearContents()
$ws.Range("O4").Value = "CAMPO VUOTO"
$ws.Range("P4").Value = "CAMPO VUOTO"

$ws.Range("A5").Value = 252334
$ws.Range("H5").ClearContents()
$ws.Range("L5").ClearContents()
$ws.Range("O5").Value = "CAMPO VUOTO"
$ws.Range("P5").Value = "CAMPO VUOTO"

$ws.Range("A6").Value = 252157
$ws.Range("H6").ClearContents()
$ws.Range("L6").ClearContents()
$ws.Range("O6").Value = "CAMPO VUOTO"
$ws.Range("P6").Value = "CAMPO VUOTO"

$ws.Range("A7").Value = 252686
$ws.Range("H7").ClearContents()
$ws.Range("L7").ClearContents()
$ws.Range("O7").Value = "CAMPO VUOTO"
$ws.Range("P7").Value = "CAMPO VUOTO"

$ws.Range("A8").Value = 252999
$ws.Range("H8").ClearContents()
$ws.Range("L8").ClearContents()
$ws.Range("O8").Value = "CAMPO VUOTO"
$ws.Range("P8").Value = "CAMPO VUOTO"

$ws.Range("A9").Value = 251231
$ws.Range("H9").Value = "CAMPO VUOTO"
$ws.Range("L9").ClearContents()
$ws.Range("O9").ClearContents()
$ws.Range("P9").ClearContents()

$ws.Range("A10").Value = 243335
$ws.Range("H10").ClearContents()
$ws.Range("L10").ClearContents()
$ws.Range("O10").Value = "CAMPO VUOTO"
$ws.Range("P10").Value = "CAMPO VUOTO"

$ws.Range("A11").Value = 251790
$ws.Range("H11").ClearContents()
$ws.Range("L11").ClearContents()
$ws.Range("O11").Value = "CAMPO VUOTO"
$ws.Range("P11").Value = "CAMPO VUOTO"

$ws.Range("A12").Value = 245623
$ws.Range("H12").Value = "CAMPO VUOTO"
$ws.Range("L12").ClearContents()
$ws.Range("O12").Value = "CAMPO VUOTO"
$ws.Range("P12").Value = "CAMPO VUOTO"

$ws.Range("A13").Value = 252652
$ws.Range("H13").Value = "CAMPO VUOTO"
$ws.Range("L13").ClearContents()
$ws.Range("O13").Value = "CAMPO VUOTO"
$ws.Range("P13").Value = "CAMPO VUOTO"

$ws.Range("A14").Value = 252899
$ws.Range("H14").ClearContents()
$ws.Range("L14").ClearContents()
$ws.Range("O14").Value = "CAMPO VUOTO"
$ws.Range("P14").Value = "CAMPO VUOTO"

$ws.Range("A15").Value = 252466
$ws.Range("H15").ClearContents()
$ws.Range("L15").ClearContents()
$ws.Range("O15").Value = "CAMPO VUOTO"
$ws.Range("P15").Value = "CAMPO VUOTO"

$ws.Range("A16").Value = 252723
$ws.Range("H16").ClearContents()
$ws.Range("L16").ClearContents()
$ws.Range("O16").Value = "CAMPO VUOTO"
$ws.Range("P16").Value = "CAMPO VUOTO"

$ws.Range("A17").Value = 251674
$ws.Range("H17").ClearContents()
$ws.Range("L17").ClearContents()
$ws.Range("O17").Value = "CAMPO VUOTO"
$ws.Range("P17").Value = "CAMPO VUOTO"

$ws.Range("A18").Value = 250284
$ws.Range("H18").ClearContents()
$ws.Range("L18").ClearContents()
$ws.Range("O18").Value = "CAMPO VUOTO"
$ws.Range("P18").Value = "CAMPO VUOTO"

$ws.Range("A19").Value = 252277
$ws.Range("H19").Value = "CAMPO VUOTO"
$ws.Range("L19").ClearContents()
$ws.Range("O19").ClearContents()
$ws.Range("P19").ClearContents()

$ws.Range("A20").Value = 252456
$ws.Range("H20").ClearContents()
$ws.Range("L20").ClearContents()
$ws.Range("O20").Value = "CAMPO VUOTO"
$ws.Range("P20").Value = "CAMPO VUOTO"

$ws.Range("A21").Value = 251919
$ws.Range("H21").ClearContents()
$ws.Range("L21").ClearContents()
$ws.Range("O21").Value = "CAMPO VUOTO"
$ws.Range("P21").Value = "CAMPO VUOTO"

$ws.Range("A22").Value = 252140
$ws.Range("H22").Value = "CAMPO VUOTO"
$ws.Range("L22").Value = "CAMPO VUOTO"
$ws.Range("O22").ClearContents()
$ws.Range("P22").ClearContents()

$ws.Range("A23").Value = 243525
$ws.Range("H23").ClearContents()
$ws.Range("L23").ClearContents()
$ws.Range("O23").Value = "CAMPO VUOTO"
$ws.Range("P23").Value = "CAMPO VUOTO"

$ws.Range("A24").Value = 243529
$ws.Range("H24").ClearContents()
$ws.Range("L24").ClearContents()
$ws.Range("O24").Value = "CAMPO VUOTO"
$ws.Range("P24").Value = "CAMPO VUOTO"

$ws.Range("A25").Value = 251798
$ws.Range("H25").ClearContents()
$ws.Range("L25").ClearContents()
$ws.Range("O25").Value = "CAMPO VUOTO"
$ws.Range("P25").Value = "CAMPO VUOTO"

$ws.Range("A26").Value = 252336
$ws.Range("H26").ClearContents()
$ws.Range("L26").ClearContents()
$ws.Range("O26").Value = "CAMPO VUOTO"
$ws.Range("P26").Value = "CAMPO VUOTO"

$ws.Range("A27").Value = 252755
$ws.Range("H27").ClearContents()
$ws.Range("L27").ClearContents()
$ws.Range("O27").Value = "CAMPO VUOTO"
$ws.Range("P27").Value = "CAMPO VUOTO"

$ws.Range("A28").Value = 252244
$ws.Range("H28").ClearContents()
$ws.Range("L28").ClearContents()
$ws.Range("O28").Value = "CAMPO VUOTO"
$ws.Range("P28").Value = "CAMPO VUOTO"

$ws.Range("A29").Value = 251685
$ws.Range("H29").ClearContents()
$ws.Range("L29").ClearContents()
$ws.Range("O29").Value = "CAMPO VUOTO"
$ws.Range("P29").Value = "CAMPO VUOTO"

$ws.Range("A30").Value = 252243
$ws.Range("H30").ClearContents()
$ws.Range("L30").ClearContents()
$ws.Range("O30").Value = "CAMPO VUOTO"
$ws.Range("P30").Value = "CAMPO VUOTO"

$ws.Range("A31").Value = 252364
$ws.Range("H31").Value = "CAMPO VUOTO"
$ws.Range("L31").ClearContents()
$ws.Range("O31").ClearContents()
$ws.Range("P31").ClearContents()

$ws.Range("A32").Value = 252549
$ws.Range("H32").ClearContents()
$ws.Range("L32").ClearContents()
$ws.Range("O32").Value = "CAMPO VUOTO"
$ws.Range("P32").Value = "CAMPO VUOTO"

$ws.Range("A33").Value = 252347
$ws.Range("H33").ClearContents()
$ws.Range("L33").ClearContents()
$ws.Range("O33").Value = "CAMPO VUOTO"
$ws.Range("P33").Value = "CAMPO VUOTO"

$ws.Range("A34").Value = 252237
$ws.Range("H34").ClearContents()
$ws.Range("L34").ClearContents()
$ws.Range("O34").Value = "CAMPO VUOTO"
$ws.Range("P34").Value = "CAMPO VUOTO"

$ws.Range("A35").Value = 252547
$ws.Range("H35").ClearContents()
$ws.Range("L35").ClearContents()
$ws.Range("O35").Value = "CAMPO VUOTO"
$ws.Range("P35").Value = "CAMPO VUOTO"

$ws.Range("A36").Value = 252467
$ws.Range("H36").ClearContents()
$ws.Range("L36").ClearContents()
$ws.Range("O36").Value = "CAMPO VUOTO"
$ws.Range("P36").Value = "CAMPO VUOTO"

$ws.Range("A37").Value = 252220
$ws.Range("H37").ClearContents()
$ws.Range("L37").ClearContents()
$ws.Range("O37").Value = "CAMPO VUOTO"
$ws.Range("P37").Value = "CAMPO VUOTO"

$ws.Range("A38").Value = 252063
$ws.Range("H38").ClearContents()
$ws.Range("L38").ClearContents()
$ws.Range("O38").Value = "CAMPO VUOTO"
$ws.Range("P38").Value = "CAMPO VUOTO"

$ws.Range("A39").Value = 251926
$ws.Range("H39").ClearContents()
$ws.Range("L39").ClearContents()
$ws.Range("O39").Value = "CAMPO VUOTO"
$ws.Range("P39").Value = "CAMPO VUOTO"

$ws.Range("A40").Value = 252418
$ws.Range("H40").Value = "CAMPO VUOTO"
$ws.Range("L40").ClearContents()
$ws.Range("O40").ClearContents()
$ws.Range("P40").ClearContents()

$ws.Range("A41").Value = 252476
$ws.Range("H41").Value = "CAMPO VUOTO"
$ws.Range("L41").ClearContents()
$ws.Range("O41").Value = "CAMPO VUOTO"
$ws.Range("P41").Value = "CAMPO VUOTO"

$ws.Range("A42").Value = 251495
$ws.Range("H42").ClearContents()
$ws.Range("L42").ClearContents()
$ws.Range("O42").Value = "CAMPO VUOTO"
$ws.Range("P42").Value = "CAMPO VUOTO"

$ws.Range("A43").Value = 251310
$ws.Range("H43").ClearContents()
$ws.Range("L43").Value = "CAMPO VUOTO"
$ws.Range("O43").Value = "CAMPO VUOTO"
$ws.Range("P43").Value = "CAMPO VUOTO"

$ws.Range("A44").Value = 252047
$ws.Range("H44").Value = "CAMPO VUOTO"
$ws.Range("L44").ClearContents()
$ws.Range("O44").Value = "CAMPO VUOTO"
$ws.Range("P44").Value = "CAMPO VUOTO"

$ws.Range("A45").Value = 252298
$ws.Range("H45").ClearContents()
$ws.Range("L45").ClearContents()
$ws.Range("O45").Value = "CAMPO VUOTO"
$ws.Range("P45").Value = "CAMPO VUOTO"

$ws.Range("A46").Value = 252146
$ws.Range("H46").Value = "CAMPO VUOTO"
$ws.Range("L46").ClearContents()
$ws.Range("O46").Value = "CAMPO VUOTO"
$ws.Range("P46").Value = "CAMPO VUOTO"

$ws.Range("A47").Value = 243524
$ws.Range("H47").ClearContents()
$ws.Range("L47").ClearContents()
$ws.Range("O47").Value = "CAMPO VUOTO"
$ws.Range("P47").Value = "CAMPO VUOTO"

$ws.Range("A48").Value = 252230
$ws.Range("H48").ClearContents()
$ws.Range("L48").ClearContents()
$ws.Range("O48").Value = "CAMPO VUOTO"
$ws.Range("P48").Value = "CAMPO VUOTO"

$ws.Range("A49").Value = 252470
$ws.Range("H49").Value = "CAMPO VUOTO"
$ws.Range("L49").ClearContents()
$ws.Range("O49").Value = "CAMPO VUOTO"
$ws.Range("P49").Value = "CAMPO VUOTO"

$ws.Range("A50").Value = 252638
$ws.Range("H50").ClearContents()
$ws.Range("L50").ClearContents()
$ws.Range("O50").Value = "CAMPO VUOTO"
$ws.Range("P50").Value = "CAMPO VUOTO"

$ws.Range("A51").Value = 252785
$ws.Range("H51").ClearContents()
$ws.Range("L51").ClearContents()
$ws.Range("O51").Value = "CAMPO VUOTO"
$ws.Range("P51").Value = "CAMPO VUOTO"

$ws.Range("A52").Value = 252546
$ws.Range("H52").ClearContents()
$ws.Range("L52").ClearContents()
$ws.Range("O52").Value = "CAMPO VUOTO"
$ws.Range("P52").Value = "CAMPO VUOTO"

$ws.Range("A53").Value = 252144
$ws.Range("H53").Value = "CAMPO VUOTO"
$ws.Range("L53").ClearContents()
$ws.Range("O53").Value = "CAMPO VUOTO"
$ws.Range("P53").Value = "CAMPO VUOTO"

$ws.Range("A54").Value = 252350
$ws.Range("H54").ClearContents()
$ws.Range("L54").ClearContents()
$ws.Range("O54").Value = "CAMPO VUOTO"
$ws.Range("P54").Value = "CAMPO VUOTO"

$ws.Range("A55").Value = 252207
$ws.Range("H55").ClearContents()
$ws.Range("L55").ClearContents()
$ws.Range("O55").Value = "CAMPO VUOTO"
$ws.Range("P55").Value = "CAMPO VUOTO"

$ws.Range("A56").Value = 252286
$ws.Range("H56").ClearContents()
$ws.Range("L56").ClearContents()
$ws.Range("O56").ClearContents()
$ws.Range("P56").ClearContents()

$ws.Range("A57").Value = 252284
$ws.Range("H57").ClearContents()
$ws.Range("L57").ClearContents()
$ws.Range("O57").ClearContents()
$ws.Range("P57").ClearContents()

$ws.Range("A58").Value = 252814
$ws.Range("H58").ClearContents()
$ws.Range("L58").ClearContents()
$ws.Range("O58").Value = "CAMPO VUOTO"
$ws.Range("P58").Value = "CAMPO VUOTO"

$ws.Range("A59").Value = 252371
$ws.Range("H59").ClearContents()
$ws.Range("L59").ClearContents()
$ws.Range("O59").Value = "CAMPO VUOTO"
$ws.Range("P59").Value = "CAMPO VUOTO"

$ws.Range("A60").Value = 252665
$ws.Range("H60").Value = "CAMPO VUOTO"
$ws.Range("L60").ClearContents()
$ws.Range("O60").Value = "CAMPO VUOTO"
$ws.Range("P60").Value = "CAMPO VUOTO"

$ws.Range("A61").Value = 252285
$ws.Range("H61").ClearContents()
$ws.Range("L61").ClearContents()
$ws.Range("O61").ClearContents()
$ws.Range("P61").ClearContents()

$ws.Range("A62").Value = 252071
$ws.Range("H62").ClearContents()
$ws.Range("L62").ClearContents()
$ws.Range("O62").Value = "CAMPO VUOTO"
$ws.Range("P62").Value = "CAMPO VUOTO"

$ws.Range("A63").Value = 244828
$ws.Range("H63").ClearContents()
$ws.Range("L63").ClearContents()
$ws.Range("O63").Value = "CAMPO VUOTO"
$ws.Range("P63").Value = "CAMPO VUOTO"

$ws.Range("A64").Value = 241783
$ws.Range("H64").ClearContents()
$ws.Range("L64").ClearContents()
$ws.Range("O64").Value = "CAMPO VUOTO"
$ws.Range("P64").Value = "CAMPO VUOTO"

$ws.Range("A65").Value = 252345
$ws.Range("H65").ClearContents()
$ws.Range("L65").ClearContents()
$ws.Range("O65").Value = "CAMPO VUOTO"
$ws.Range("P65").Value = "CAMPO VUOTO"

$ws.Range("A66").Value = 252997
$ws.Range("H66").ClearContents()
$ws.Range("L66").ClearContents()
$ws.Range("O66").Value = "CAMPO VUOTO"
$ws.Range("P66").Value = "CAMPO VUOTO"

$ws.Range("A67").Value = 252417
$ws.Range("H67").Value = "CAMPO VUOTO"
$ws.Range("L67").ClearContents()
$ws.Range("O67").Value = "CAMPO VUOTO"
$ws.Range("P67").Value = "CAMPO VUOTO"

$ws.Range("A68").Value = 252087
$ws.Range("H68").ClearContents()
$ws.Range("L68").ClearContents()
$ws.Range("O68").Value = "CAMPO VUOTO"
$ws.Range("P68").Value = "CAMPO VUOTO"

$ws.Range("A69").Value = 252201
$ws.Range("H69").ClearContents()
$ws.Range("L69").ClearContents()
$ws.Range("O69").Value = "CAMPO VUOTO"
$ws.Range("P69").Value = "CAMPO VUOTO"

$ws.Range("A70").Value = 251849
$ws.Range("H70").Value = "CAMPO VUOTO"
$ws.Range("L70").ClearContents()
$ws.Range("O70").ClearContents()
$ws.Range("P70").ClearContents()

$ws.Range("A71").Value = 252085
$ws.Range("H71").ClearContents()
$ws.Range("L71").ClearContents()
$ws.Range("O71").Value = "CAMPO VUOTO"
$ws.Range("P71").Value = "CAMPO VUOTO"

$ws.Range("A72").Value = 252152
$ws.Range("H72").ClearContents()
$ws.Range("L72").ClearContents()
$ws.Range("O72").Value = "CAMPO VUOTO"
$ws.Range("P72").Value = "CAMPO VUOTO"

$ws.Range("A73").Value = 252784
$ws.Range("H73").ClearContents()
$ws.Range("L73").ClearContents()
$ws.Range("O73").Value = "CAMPO VUOTO"
$ws.Range("P73").Value = "CAMPO VUOTO"

$ws.Range("A74").Value = 252402
$ws.Range("H74").ClearContents()
$ws.Range("L74").ClearContents()
$ws.Range("O74").Value = "CAMPO VUOTO"
$ws.Range("P74").Value = "CAMPO VUOTO"

$ws.Range("A75").Value = 252978
$ws.Range("H75").ClearContents()
$ws.Range("L75").ClearContents()
$ws.Range("O75").Value = "CAMPO VUOTO"
$ws.Range("P75").Value = "CAMPO VUOTO"

$ws.Range("A76").Value = 243569
$ws.Range("H76").ClearContents()
$ws.Range("L76").ClearContents()
$ws.Range("O76").Value = "CAMPO VUOTO"
$ws.Range("P76").Value = "CAMPO VUOTO"

$ws.Range("A77").Value = 252980
$ws.Range("H77").ClearContents()
$ws.Range("L77").ClearContents()
$ws.Range("O77").Value = "CAMPO VUOTO"
$ws.Range("P77").Value = "CAMPO VUOTO"

$ws.Range("A78").Value = 252983
$ws.Range("H78").ClearContents()
$ws.Range("L78").Value = "CAMPO VUOTO"
$ws.Range("O78").Value = "CAMPO VUOTO"
$ws.Range("P78").Value = "CAMPO VUOTO"

$ws.Range("A79").Value = 252790
$ws.Range("H79").ClearContents()
$ws.Range("L79").ClearContents()
$ws.Range("O79").Value = "CAMPO VUOTO"
$ws.Range("P79").Value = "CAMPO VUOTO"

$ws.Range("A80").Value = 252783
$ws.Range("H80").ClearContents()
$ws.Range("L80").ClearContents()
$ws.Range("O80").Value = "CAMPO VUOTO"
$ws.Range("P80").Value = "CAMPO VUOTO"

$ws.Range("A81").Value = 252636
$ws.Range("H81").ClearContents()
$ws.Range("L81").ClearContents()
$ws.Range("O81").Value = "CAMPO VUOTO"
$ws.Range("P81").Value = "CAMPO VUOTO"

$ws.Range("A82").Value = 244023
$ws.Range("H82").ClearContents()
$ws.Range("L82").ClearContents()
$ws.Range("O82").Value = "CAMPO VUOTO"
$ws.Range("P82").Value = "CAMPO VUOTO"

